# Alteração para diferentes tipos de bimestrais
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 18 (20:00)
$ws.Range("B18").Value = "[Guilherme-C. L. P.-2NB, Guilherme-C. L. P.-2NB, Guilherme-C. L. P.-2NB, Guilherme-C. L. P.-2NB]"
$ws.Range("C18").Value = "Suzanny-Des. Maq. CAD-"
$ws.Range("D18").Value = "[Aderci-Fresagem-2NB, Aderci-Fresagem-2NB, Aderci-Fresagem-2NB, Aderci-Fresagem-2NB]"
$ws.Range("E18").Value = "[Ismail-Metrologia 2-2NB, Ismail-Metrologia 2-2NB, Ismail-Metrologia 2-2NB, Ismail-Metrologia 2-2NB]"
$ws.Range("F18").Value = "Claudinei-Elemaq.-"

# Row 19 (20:50)
$ws.Range("B19").Value = "[Leandro-M.S.R.A.C.-2NB, Leandro-M.S.R.A.C.-2NB, Leandro-M.S.R.A.C.-2NB, Leandro-M.S.R.A.C.-2NB]"
$ws.Range("C19").Value = "Suzanny-Des. Maq. CAD-"
$ws.Range("D19").Value = "[Guilherme-Coman. Hidraulicos-2NB, Guilherme-Coman. Hidraulicos-2NB, Guilherme-Coman. Hidraulicos-2NB, Guilherme-Coman. Hidraulicos-2NB]"
$ws.Range("E19").Value = "Claudinei-Elemaq.-"
$ws.Range("F19").Value = "[Leandro-M. Maq. E. I.-2NB, Leandro-M. Maq. E. I.-2NB, Leandro-M. Maq. E. I.-2NB, Leandro-M. Maq. E. I.-2NB]"

# Row 20 (21:40)
$ws.Range("B20").Value = "[Victor S.-Usin. CNC-2NB, Victor S.-Usin. CNC-2NB, Victor S.-Usin. CNC-2NB, Victor S.-Usin. CNC-2NB]"
$ws.Range("C20").Value = "[Paulo Rob.-CAD / CAM-2NB, Paulo Rob.-CAD / CAM-2NB, Paulo Rob.-CAD / CAM-2NB, Paulo Rob.-CAD / CAM-2NB]"
$ws.Range("D20").Value = "[Rogério-Retífica-2NB, Rogério-Retífica-2NB, Rogério-Retífica-2NB, Rogério-Retífica-2NB]"
$ws.Range("E20").Value = "Euclides-Gest. Int.-"
$ws.Range("F20").Value = "Suzanny-Des. Maq. CAD-"

# Row 21 (22:35)
$ws.Range("C21").Value = "[Elcio D.-C. Pneumática-2NB, Elcio D.-C. Pneumática-2NB, Elcio D.-C. Pneumática-2NB, Elcio D.-C. Pneumática-2NB]"
$ws.Range("D21").Value = "[Joel L.-Fundição-2NB, Joel L.-Fundição-2NB, Joel L.-Fundição-2NB, Joel L.-Fundição-2NB]"
$ws.Range("E21").Value = "Euclides-Gest. Int.-"
$ws.Range("F21").Value = "[Cláudio-Soldagem-2NB, Cláudio-Soldagem-2NB, Cláudio-Soldagem-2NB, Cláudio-Soldagem-2NB]"
